$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was updated
# from 45190 (2023-09-21) to 45192 (2023-09-23) for every data row
# (rows 2 through 310).
$newValue = 45192
$lastRow = 310

for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 3).Value = $newValue
}
